$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-CellText $ws.Range("D2") "60.104.02"
Set-CellText $ws.Range("E2") "  +0.21%  "
Set-CellText $ws.Range("D3") "2.414.39"
Set-CellText $ws.Range("E3") "  -0.22%  "
Set-CellText $ws.Range("E4") "  -0.03%  "
Set-CellText $ws.Range("D5") "554.19"
Set-CellText $ws.Range("E5") "  +0.37%  "
Set-CellText $ws.Range("D6") "136.81"
Set-CellText $ws.Range("E6") "  -0.32%  "
Set-CellText $ws.Range("E7") "  +0.02%  "
Set-CellText $ws.Range("E8") "  +0.52%  "
Set-CellText $ws.Range("E9") "  -1.01%  "
Set-CellText $ws.Range("D10") "5.67"
Set-CellText $ws.Range("E10") "  -1.87%  "
Set-CellText $ws.Range("E11") "  -0.68%  "
Set-CellText $ws.Range("E12") "  -1.20%  "
Set-CellText $ws.Range("D13") "24.76"
Set-CellText $ws.Range("E13") "  -0.12%  "
Set-CellText $ws.Range("D14") "2.846.97"
Set-CellText $ws.Range("E14") "  -0.18%  "
Set-CellText $ws.Range("D15") "59.992.10"
Set-CellText $ws.Range("E15") "  +0.10%  "
Set-CellText $ws.Range("D17") "2.414.54"
Set-CellText $ws.Range("E17") "  -0.57%  "
Set-CellText $ws.Range("E18") "  -0.81%  "
Set-CellText $ws.Range("D19") "4.53"
Set-CellText $ws.Range("E19") "  +3.44%  "
Set-CellText $ws.Range("D20") "326.80"
Set-CellText $ws.Range("E20") "  -1.35%  "
Set-CellText $ws.Range("E21") "  +1.12%  "
Set-CellText $ws.Range("D22") "1.00"
Set-CellText $ws.Range("E22") "  +0.06%  "
Set-CellText $ws.Range("D23") "64.79"
Set-CellText $ws.Range("E23") "  -1.29%  "
Set-CellText $ws.Range("E24") "  +5.88%  "
Set-CellText $ws.Range("E25") "  +0.13%  "
Set-CellText $ws.Range("E26") "  -0.01%  "
Set-CellText $ws.Range("D27") "1.41"
Set-CellText $ws.Range("E27") "  +5.03%  "
Set-CellText $ws.Range("E28") "  -1.18%  "
Set-CellText $ws.Range("D29") "1.78"
Set-CellText $ws.Range("E29") "  +0.19%  "
Set-CellText $ws.Range("D30") "170.67"
Set-CellText $ws.Range("E30") "  +0.32%  "
Set-CellText $ws.Range("E31") "  -1.74%  "
Set-CellText $ws.Range("D32") "1.08"
Set-CellText $ws.Range("E32") "  +5.27%  "
Set-CellText $ws.Range("E33") "  -3.06%  "
Set-CellText $ws.Range("E34") "  -0.74%  "
Set-CellText $ws.Range("E35") "  +0.04%  "
Set-CellText $ws.Range("E36") "  +1.97%  "
Set-CellText $ws.Range("B37") "FirstDigitalUSD"
Set-CellText $ws.Range("C37") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-CellText $ws.Range("D37") "1.00"
Set-CellText $ws.Range("E37") "  +0.02%  "
Set-CellText $ws.Range("B38") "NEARProtocol"
Set-CellText $ws.Range("C38") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-CellText $ws.Range("D38") "4.22"
Set-CellText $ws.Range("E38") "  +0.99%  "
Set-CellText $ws.Range("D39") "324.76"
Set-CellText $ws.Range("E39") "  +3.86%  "
Set-CellText $ws.Range("E40") "  -1.00%  "
Set-CellText $ws.Range("D41") "146.44"
Set-CellText $ws.Range("E41") "  +4.73%  "
Set-CellText $ws.Range("D42") "3.63"
Set-CellText $ws.Range("E42") "  -1.20%  "
Set-CellText $ws.Range("E43") "  -0.13%  "
Set-CellText $ws.Range("D44") "19.77"
Set-CellText $ws.Range("E44") "  +2.41%  "
Set-CellText $ws.Range("E45") "  -0.69%  "
Set-CellText $ws.Range("E46") "  +0.28%  "
Set-CellText $ws.Range("E47") "  -1.19%  "
Set-CellText $ws.Range("D48") "11.04"
Set-CellText $ws.Range("E48") "  -0.02%  "
Set-CellText $ws.Range("E49") "  -0.87%  "
Set-CellText $ws.Range("E50") "  -0.60%  "
Set-CellText $ws.Range("D51") "0.939"
Set-CellText $ws.Range("E51") "  -1.75%  "
